# CT - Manter Sprint : add the three "unsuccessful" test-case rows to the
# "Scripts" sheet, make that sheet the active one (tab selected) with
# B8 as the selected cell, matching the target commit.

$wb  = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item(3)   # Scripts

# --- New content for rows 5, 6 and 7 -------------------------------------

# Row 5 - "Usuário sem permissão" test case
$ws3.Range("B5").Value = "Efetuar sem sucesso o cadastro de um Sprint. Usuário sem permissão a funcionalidade."
$ws3.Range("C5").Value = "1- Tentar executar o passo #1 com usuário sem acesso a funcionalidade."
$ws3.Range("D5").Value = "O sistema não deverá exibir a opção para cadastro de Sprint."

# Row 6 - "Dados inválidos" test case
$ws3.Range("B6").Value = "Efetuar sem sucesso o cadastro de um Sprint. Dados inválidos.`nPasso #1 deverá ter sido executado com sucesso."
$ws3.Range("C6").Value = "1- Preencher os campos de cadastro da Sprint com dados inválidos.`n2- Acionar o Salvar."
$ws3.Range("D6").Value = "O sistema não deverá salvar a Sprint e deverá exibir uma mensagem de erro informando os campos com preenchimento inválido."

# Row 7 - "Campos obrigatórios não preenchidos" test case
$ws3.Range("B7").Value = "Efetuar sem sucesso o cadastro de um Sprint. Campos obrigatórios não preenchidos.`nPasso #1 deverá ter sido executado com sucesso."
$ws3.Range("C7").Value = "1- Não preencher os campos obrigatórios do cadastro da Sprint.`n2- Acionar o Salvar."
$ws3.Range("D7").Value = "O sistema não deverá salvar a Sprint e deverá exibir uma mensagem de erro informando que os campos estão em branco."

# --- Formatting: wrap + vertically centered text, taller rows ------------

$fmtRange = $ws3.Range("B5:D7")
$fmtRange.WrapText = $true
$fmtRange.VerticalAlignment = -4108   # xlCenter

$ws3.Rows.Item(5).RowHeight = 25.5
$ws3.Rows.Item(6).RowHeight = 38.25
$ws3.Rows.Item(7).RowHeight = 51

# --- Active sheet / selection --------------------------------------------
# Originally "Capa" (sheet1) was the selected tab; now it is "Scripts"
# (sheet3), with B8 selected.

$ws3.Activate()
$ws3.Range("B8").Select()
